# Daily attendance processing - 2026-01-22 07:17:51
# Normalizes the "Recorded By" (column G) author lists: the most recent
# recorder (last name in the comma-separated list) is promoted to the
# front of the list, e.g. "System, dnasr281@gmail.com" becomes
# "dnasr281@gmail.com, System".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count
$recordedByCol = 7  # column G: "Recorded By"

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $recordedByCol)
    $val = $cell.Value2

    if ($val -eq $null) { continue }
    if ($val -eq "") { continue }
    # This particular author combination is left as-is.
    if ($val -eq "backup@backdoor.com, System") { continue }

    $parts = @($val -split ", ")
    if ($parts.Count -le 1) { continue }

    $lastAuthor = $parts[$parts.Count - 1]
    $remaining = $parts[0..($parts.Count - 2)]
    $newParts = @($lastAuthor) + $remaining
    $newVal = $newParts -join ", "

    if ($newVal -ne $val) {
        $cell.Value = $newVal
    }
}

Write-Output "Recorded By normalization complete"
